# Auto-generated edit script applying numeric updates to Kujata_Profits workbook
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1730.25
$ws.Range("J29").Value = 2253.3333
$ws.Range("L29").Value = 6759.999899999999
$ws.Range("N29").Value = -7321.999899999999
$ws.Range("H38").Value = 1560.7561
$ws.Range("J38").Value = 2003.5161
$ws.Range("L38").Value = 6010.5483
$ws.Range("N38").Value = -6754.5483
$ws.Range("H43").Value = 9279501
$ws.Range("I43").Value = 50375.5
$ws.Range("J43").Value = 13894064
$ws.Range("K43").Value = 50375.5
$ws.Range("L43").Value = 13894064
$ws.Range("M43").Value = -50306.5
$ws.Range("N43").Value = -13894202
$ws.Range("H113").Value = 2099
$ws.Range("I113").Value = 1799
$ws.Range("J113").Value = 2249
$ws.Range("K113").Value = 1799
$ws.Range("L113").Value = 2249
$ws.Range("M113").Value = 1455
$ws.Range("N113").Value = -8757
$ws.Range("H115").Value = 592.5
$ws.Range("I115").Value = 456.66666
$ws.Range("J115").Value = 1000
$ws.Range("K115").Value = 1369.99998
$ws.Range("L115").Value = 3000
$ws.Range("M115").Value = 197.0000199999999
$ws.Range("N115").Value = -6134
$ws.Range("H138").Value = 1604.9193
$ws.Range("I138").Value = 1100.5454
$ws.Range("J138").Value = 1713.7059
$ws.Range("K138").Value = 3301.6362
$ws.Range("L138").Value = 5141.1177
$ws.Range("M138").Value = 1838.3638
$ws.Range("N138").Value = -15421.1177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 10000000
$ws.Range("I8").Value = 10000000
$ws.Range("K8").Value = 10000000
$ws.Range("M8").Value = -9999856
$ws.Range("H13").Value = 12500251
$ws.Range("I13").Value = 16666667
$ws.Range("J13").Value = 1004
$ws.Range("K13").Value = 16666667
$ws.Range("L13").Value = 1004
$ws.Range("M13").Value = -16666523
$ws.Range("N13").Value = -1292
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H32").Value = 4724.227
$ws.Range("I32").Value = 4787.5815
$ws.Range("K32").Value = 4787.5815
$ws.Range("M32").Value = -4500.5815
$ws.Range("H45").Value = 1324.7142
$ws.Range("I45").Value = 1184.7273
$ws.Range("K45").Value = 1184.7273
$ws.Range("M45").Value = -807.7273
$ws.Range("H61").Value = 250001120
$ws.Range("I61").Value = 500000260
$ws.Range("K61").Value = 500000260
$ws.Range("M61").Value = -500000048
$ws.Range("H88").Value = 2424
$ws.Range("I88").Value = 1800.4
$ws.Range("J88").Value = 3203.5
$ws.Range("K88").Value = 1800.4
$ws.Range("L88").Value = 3203.5
$ws.Range("M88").Value = -1394.4
$ws.Range("N88").Value = -4015.5
$ws.Range("H91").Value = 2424
$ws.Range("I91").Value = 1800.4
$ws.Range("J91").Value = 3203.5
$ws.Range("K91").Value = 1800.4
$ws.Range("L91").Value = 3203.5
$ws.Range("M91").Value = -396.4000000000001
$ws.Range("N91").Value = -6011.5
$ws.Range("H97").Value = 621.6667
$ws.Range("I97").Value = 471.15384
$ws.Range("K97").Value = 471.15384
$ws.Range("M97").Value = 24.84616
$ws.Range("H110").Value = 2713.375
$ws.Range("I110").Value = 594
$ws.Range("J110").Value = 3985
$ws.Range("K110").Value = 594
$ws.Range("L110").Value = 3985
$ws.Range("M110").Value = 1451
$ws.Range("N110").Value = -8075
$ws.Range("H122").Value = 1874.0769
$ws.Range("I122").Value = 1304.4546
$ws.Range("K122").Value = 3913.3638
$ws.Range("M122").Value = -1463.3638
$ws.Range("H132").Value = 3074.3928
$ws.Range("I132").Value = 2845.158
$ws.Range("K132").Value = 8535.474
$ws.Range("M132").Value = -6005.474
$ws.Range("H136").Value = 250001120
$ws.Range("I136").Value = 500000260
$ws.Range("K136").Value = 1500000780
$ws.Range("M136").Value = -1499998230

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = 0
$ws.Range("H25").Value = 650
$ws.Range("I25").Value = 650
$ws.Range("K25").Value = 650
$ws.Range("M25").Value = -415
$ws.Range("H29").Value = 744
$ws.Range("I29").Value = 744
$ws.Range("K29").Value = 744
$ws.Range("M29").Value = -455
$ws.Range("H99").Value = 166667820
$ws.Range("I99").Value = 333334340
$ws.Range("K99").Value = 333334340
$ws.Range("M99").Value = -333332842
$ws.Range("H134").Value = 11319.091
$ws.Range("I134").Value = 1362.4
$ws.Range("J134").Value = 19616.334
$ws.Range("K134").Value = 4087.2
$ws.Range("L134").Value = 58849.00199999999
$ws.Range("M134").Value = -1552.2
$ws.Range("N134").Value = -63919.00199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 470
$ws.Range("I105").Value = 470
$ws.Range("K105").Value = 470
$ws.Range("M105").Value = 1277

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 58823670
$ws.Range("J12").Value = 112.833336
$ws.Range("L12").Value = 338.500008
$ws.Range("N12").Value = -684.500008
$ws.Range("H23").Value = 311.5
$ws.Range("J23").Value = 367.5
$ws.Range("L23").Value = 1102.5
$ws.Range("N23").Value = -1572.5
$ws.Range("H39").Value = 2638.625
$ws.Range("J39").Value = 2533.0454
$ws.Range("L39").Value = 7599.1362
$ws.Range("N39").Value = -8187.1362
$ws.Range("H55").Value = 3183.3333
$ws.Range("J55").Value = 3183.3333
$ws.Range("L55").Value = 9549.999899999999
$ws.Range("N55").Value = -9903.999899999999
$ws.Range("H92").Value = 239.65517
$ws.Range("I92").Value = 230.8
$ws.Range("J92").Value = 295
$ws.Range("K92").Value = 692.4000000000001
$ws.Range("L92").Value = 885
$ws.Range("M92").Value = 555.5999999999999
$ws.Range("N92").Value = -3381
$ws.Range("H131").Value = 18185182
$ws.Range("J131").Value = 3879.0852
$ws.Range("L131").Value = 11637.2556
$ws.Range("N131").Value = -21717.2556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 240.54546
$ws.Range("I2").Value = 185.57143
$ws.Range("J2").Value = 336.75
$ws.Range("K2").Value = 185.57143
$ws.Range("L2").Value = 336.75
$ws.Range("M2").Value = -72.57142999999999
$ws.Range("N2").Value = -562.75
$ws.Range("H20").Value = 32500000
$ws.Range("I20").Value = 32500000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 32500000
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -32499755
$ws.Range("H24").Value = 17501500
$ws.Range("I24").Value = 23333334
$ws.Range("J24").Value = 6000
$ws.Range("K24").Value = 23333334
$ws.Range("L24").Value = 6000
$ws.Range("M24").Value = -23333161
$ws.Range("N24").Value = -6346
$ws.Range("H122").Value = 1272.675
$ws.Range("I122").Value = 1205.9667
$ws.Range("J122").Value = 1472.8
$ws.Range("K122").Value = 3617.9001
$ws.Range("L122").Value = 4418.4
$ws.Range("M122").Value = -1167.9001
$ws.Range("N122").Value = -9318.4
$ws.Range("H132").Value = 6302.875
$ws.Range("I132").Value = 8352.75
$ws.Range("K132").Value = 25058.25
$ws.Range("M132").Value = -22528.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 638.1667
$ws.Range("I22").Value = 432
$ws.Range("J22").Value = 844.3333
$ws.Range("K22").Value = 432
$ws.Range("L22").Value = 844.3333
$ws.Range("M22").Value = -137
$ws.Range("N22").Value = -1434.3333
$ws.Range("H27").Value = 638.1667
$ws.Range("I27").Value = 432
$ws.Range("J27").Value = 844.3333
$ws.Range("K27").Value = 432
$ws.Range("L27").Value = 844.3333
$ws.Range("M27").Value = -325
$ws.Range("N27").Value = -1058.3333
$ws.Range("H68").Value = 1713.1904
$ws.Range("I68").Value = 1705.1052
$ws.Range("J68").Value = 1790
$ws.Range("K68").Value = 1705.1052
$ws.Range("L68").Value = 1790
$ws.Range("M68").Value = -956.1052
$ws.Range("N68").Value = -3288
$ws.Range("H71").Value = 1713.1904
$ws.Range("I71").Value = 1705.1052
$ws.Range("J71").Value = 1790
$ws.Range("K71").Value = 8525.526
$ws.Range("L71").Value = 8950
$ws.Range("M71").Value = -4781.526
$ws.Range("N71").Value = -16438
$ws.Range("H93").Value = 1066.6666
$ws.Range("H100").Value = 1389
$ws.Range("I100").Value = 1389
$ws.Range("K100").Value = 1389
$ws.Range("M100").Value = -848

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7000
$ws.Range("J15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("N15").Value = -7576
$ws.Range("H18").Value = 2000
$ws.Range("I18").Value = 2000
$ws.Range("K18").Value = 2000
$ws.Range("M18").Value = -1827
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").ClearContents()
$ws.Range("N31").Value = 0
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0
$ws.Range("H106").Value = 16500
$ws.Range("J106").Value = 16500
$ws.Range("L106").Value = 16500
$ws.Range("M106").Value = -19024
$ws.Range("H109").Value = 36650.8
$ws.Range("J109").Value = 33228
$ws.Range("L109").Value = 33228
$ws.Range("N109").Value = -36002
$ws.Range("H132").Value = 8573.637000000001
$ws.Range("I132").Value = 14902.4
$ws.Range("K132").Value = 44707.2
$ws.Range("M132").Value = -42177.2
